$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 previously mirrored A1's text ("DPG458") using a quote-prefixed text
# style (s="1"). Update it to a new text value "Amol" while preserving
# that quote-prefix/text style - the leading apostrophe forces Excel to
# treat the value as text with the quote-prefix style retained.
$ws.Range("D1").Value = "'Amol"
